# Generate Report for Handoff
#
# The c4e3a62b-d9fe-4fda-8852-3a931081d1e1 file (row 3 on every sheet) has
# moved from "Handed back: in sync with en-US" to "Ready for handoff", and
# its handoff timestamps were refreshed accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 = c4e3a62b-....md ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 12:12:47"

# --- zh-cn sheet: row 3 = c4e3a62b-....md ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 12:12:39"

# --- de-de sheet: row 3 = c4e3a62b-....md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 12:12:47"
